$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reformat the date-time literals embedded in the migration strings ---
# "...T09:00:00Z"  ->  "... 09:00:00"
# "...T17:45:00Z"  ->  "... 17:45:00"
# Use Range.Replace so the shared-string table entries are edited in place
# (preserves ordering / does not create brand-new shared-string records).
$dataRange = $ws.Range("A1:A16")
$dataRange.Replace("T09:00:00Z", " 09:00:00")
$dataRange.Replace("T17:45:00Z", " 17:45:00")

# --- 2. Update the sheet view: scroll so column B is first visible, and move
#        the active selection to K9 ---
$ws.Range("K9").Select()
